# Update Weight_matrix.xlsx: add 5 new "AI image" model columns (Y:AC)
# with header labels and weight values (many still NaN, per commit message).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new model names in columns Y..AC ---
$headerCols = @("Y","Z","AA","AB","AC")
$headerVals = @("MobileNetV3Small","RegNetY080","ConvNeXtSmall","EfficientNetB0","EfficientNetV2S")

for ($i = 0; $i -lt $headerCols.Length; $i++) {
    $cell = $ws.Range($headerCols[$i] + "1")
    $cell.Value = $headerVals[$i]
    $cell.HorizontalAlignment = -4152
}

# --- Data rows 2..41: same value repeated across the 5 new columns ---
$rowValues = @{
    2  = 1;  3  = 1;  4  = 0;  5  = 0;  6  = 0;
    7  = 1;  8  = 0;  9  = 0;  10 = 0;  11 = 0;
    12 = 0;  13 = 0;  14 = 1;  15 = "NaN"; 16 = 1;
    17 = "NaN"; 18 = "NaN"; 19 = "NaN"; 20 = 1; 21 = 1;
    22 = 1;  23 = 0;  24 = 0;  25 = 0;  26 = "NaN";
    27 = "NaN"; 28 = "NaN"; 29 = 1; 30 = 1; 31 = 1;
    32 = 1;  33 = 1;  34 = 1;  35 = 1;  36 = 1;
    37 = 1;  38 = 1;  39 = 1;  40 = 1;  41 = 1
}

for ($r = 2; $r -le 41; $r++) {
    $val = $rowValues[$r]
    for ($i = 0; $i -lt $headerCols.Length; $i++) {
        $cell = $ws.Range($headerCols[$i] + $r)
        $cell.Value = $val
        $cell.HorizontalAlignment = -4152
    }
}

# --- View/window cosmetics captured in the diff ---
$ws.Range("S6").Select()
$excel.ActiveWindow.Zoom = 85
